# Update cryptocurrency price/volume data (and row 47/48 coin swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data cells in columns B:E are stored as text in this sheet (coin names,
# links, and price/volume strings that often look numeric, e.g. "27.146.92" or
# "2.280"). Force text number-format before writing so Excel does not silently
# reinterpret/round them as numbers and drop significant trailing zeros.
$cells = @()
$cells += "D2"
$cells += "E2"
$cells += "D3"
$cells += "E3"
$cells += "D4"
$cells += "E4"
$cells += "D5"
$cells += "E5"
$cells += "D6"
$cells += "E6"
$cells += "E7"
$cells += "D8"
$cells += "E8"
$cells += "D9"
$cells += "E9"
$cells += "D10"
$cells += "E10"
$cells += "D11"
$cells += "E11"
$cells += "D12"
$cells += "E12"
$cells += "D13"
$cells += "E13"
$cells += "D14"
$cells += "E14"
$cells += "D15"
$cells += "E15"
$cells += "D16"
$cells += "E16"
$cells += "D17"
$cells += "E17"
$cells += "D18"
$cells += "E18"
$cells += "D19"
$cells += "E19"
$cells += "D20"
$cells += "E20"
$cells += "D21"
$cells += "E21"
$cells += "D22"
$cells += "E22"
$cells += "E23"
$cells += "D24"
$cells += "E24"
$cells += "D25"
$cells += "E25"
$cells += "D26"
$cells += "E26"
$cells += "E27"
$cells += "E28"
$cells += "D29"
$cells += "E29"
$cells += "D30"
$cells += "E30"
$cells += "D31"
$cells += "E31"
$cells += "E32"
$cells += "D33"
$cells += "E33"
$cells += "D34"
$cells += "E34"
$cells += "E35"
$cells += "D36"
$cells += "E36"
$cells += "D37"
$cells += "E37"
$cells += "D38"
$cells += "E38"
$cells += "D39"
$cells += "E39"
$cells += "E40"
$cells += "E41"
$cells += "D42"
$cells += "E42"
$cells += "D43"
$cells += "E43"
$cells += "D44"
$cells += "E44"
$cells += "E45"
$cells += "E46"
$cells += "B47"
$cells += "C47"
$cells += "D47"
$cells += "E47"
$cells += "B48"
$cells += "C48"
$cells += "D48"
$cells += "E48"
$cells += "E49"
$cells += "D50"
$cells += "E50"
$cells += "D51"
$cells += "E51"

foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.146.92'
$ws.Range("E2").Value = '  +0.12%  '

$ws.Range("D3").Value = '1.902.15'
$ws.Range("E3").Value = '  +0.59%  '

$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '305.89'
$ws.Range("E5").Value = '  -0.51%  '

$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("E7").Value = '  +1.63%  '

$ws.Range("D8").Value = '0.3761'
$ws.Range("E8").Value = '  +0.79%  '

$ws.Range("D9").Value = '0.07247'
$ws.Range("E9").Value = '  +0.32%  '

$ws.Range("D10").Value = '21.11'
$ws.Range("E10").Value = '  -0.44%  '

$ws.Range("D11").Value = '0.9019'
$ws.Range("E11").Value = '  -0.57%  '

$ws.Range("D12").Value = '0.08485'
$ws.Range("E12").Value = '  +11.15%  '

$ws.Range("D13").Value = '1.922.69'
$ws.Range("E13").Value = '  +1.76%  '

$ws.Range("D14").Value = '95.02'
$ws.Range("E14").Value = '  +0.31%  '

$ws.Range("D15").Value = '5.287'
$ws.Range("E15").Value = '  +0.19%  '

$ws.Range("D16").Value = '0.9997'
$ws.Range("E16").Value = '  -0.07%  '

$ws.Range("D17").Value = '0.000008623'
$ws.Range("E17").Value = '  +1.23%  '

$ws.Range("D18").Value = '14.52'
$ws.Range("E18").Value = '  +0.30%  '

$ws.Range("D19").Value = '0.9994'
$ws.Range("E19").Value = '  +0.01%  '

$ws.Range("D20").Value = '27.190.78'
$ws.Range("E20").Value = '  +0.15%  '

$ws.Range("D21").Value = '5.063'
$ws.Range("E21").Value = '  -0.31%  '

$ws.Range("D22").Value = '2.155.42'
$ws.Range("E22").Value = '  +1.57%  '

$ws.Range("E23").Value = '  +0.24%  '

$ws.Range("D24").Value = '6.424'
$ws.Range("E24").Value = '  -0.06%  '

$ws.Range("D25").Value = '147.29'
$ws.Range("E25").Value = '  +0.51%  '

$ws.Range("D26").Value = '2.280'
$ws.Range("E26").Value = '  +3.11%  '

$ws.Range("E27").Value = '  -2.38%  '

$ws.Range("E28").Value = '  +0.58%  '

$ws.Range("D29").Value = '114.82'
$ws.Range("E29").Value = '  +0.13%  '

$ws.Range("D30").Value = '4.810'
$ws.Range("E30").Value = '  -0.96%  '

$ws.Range("D31").Value = '4.890'
$ws.Range("E31").Value = '  -1.69%  '

$ws.Range("E32").Value = '  +0.51%  '

$ws.Range("D33").Value = '0.8081'
$ws.Range("E33").Value = '  +4.93%  '

$ws.Range("D34").Value = '0.05053'
$ws.Range("E34").Value = '  -0.58%  '

$ws.Range("E35").Value = '  -0.32%  '

$ws.Range("D36").Value = '3.447'
$ws.Range("E36").Value = '  +4.56%  '

$ws.Range("D37").Value = '2.941'
$ws.Range("E37").Value = '  -1.47%  '

$ws.Range("D38").Value = '2.626'
$ws.Range("E38").Value = '  -0.17%  '

$ws.Range("D39").Value = '0.5713'
$ws.Range("E39").Value = '  +1.61%  '

$ws.Range("E40").Value = '  -0.30%  '

$ws.Range("E41").Value = '  -0.16%  '

$ws.Range("D42").Value = '9.018'
$ws.Range("E42").Value = '  +0.58%  '

$ws.Range("D43").Value = '6.625'
$ws.Range("E43").Value = '  -0.31%  '

$ws.Range("D44").Value = '116.38'
$ws.Range("E44").Value = '  -1.37%  '

$ws.Range("E45").Value = '  +0.15%  '

$ws.Range("E46").Value = '  +0.76%  '

$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '0.9994'
$ws.Range("E47").Value = '  +0.07%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '10.11'
$ws.Range("E48").Value = '  -1.24%  '

$ws.Range("E49").Value = '  +1.09%  '

$ws.Range("D50").Value = '37.48'
$ws.Range("E50").Value = '  +0.16%  '

$ws.Range("D51").Value = '63.93'
$ws.Range("E51").Value = '  -0.67%  '
